$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A10").Value = "Running Original AI2 QA code again for comparison purposes `nbatchsize = 60"
$ws.Range("C10").Value = "DLT1 / 6"
$ws.Range("D10").Value = "EXP8.txt"
$ws.Range("B10").Value = "python -m basic.cli --mode train --noload --len_opt --cluster --batch_size 60 --run_id 5 |& tee /home/hpalangi/QA/TPR_Stuff/Codes/TPR_ver1.0/Log_Files/EXP8.txt"
$ws.Range("E10").Value = 5
$ws.Range("F10").Value = 6

$ws.Range("A10").WrapText = $true
$ws.Rows.Item(10).RowHeight = 60

$ws.Range("C14").Select()
